# Update investment-cost results values as provided by the server re-run.
# Only numeric result cells in row 2 of sheets "2025", "2030", "2035",
# "2040" and "2045" change; sheet "2050" is untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "B2" = 9739.537847600008
        "E2" = 289823.7596598056
        "I2" = 161752.8135478
        "L2" = 485245.29503538
        "M2" = 105905.87968015
        "N2" = 70831.955579581
        "O2" = 69610.4422391004
    }
    "2030" = @{
        "B2" = 47386.06393082884
        "E2" = 271236.7992183856
        "I2" = 280426.171173861
        "L2" = 184420.4799505123
        "M2" = 113936.92264746
        "N2" = 33931.8246116005
        "O2" = 50485.47232467777
    }
    "2035" = @{
        "A2" = 28619.61401238371
        "B2" = 23143.29485244409
        "E2" = 111916.8406725409
        "I2" = 150385.2728707001
        "M2" = 34803.41203795493
        "N2" = 44938.11408779013
        "O2" = 26938.31306104351
    }
    "2040" = @{
        "N2" = 1014.766490779938
    }
    "2045" = @{
        "A2" = 34409.11717595647
        "N2" = 5182.698656944208
        "O2" = 22972.54525065906
    }
}

foreach ($sheetName in $updates.Keys) {
    # Hashtable keys that look numeric (e.g. "2025") get coerced to Int64,
    # so force back to string before using as a worksheet name lookup.
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
